$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 98698
$ws.Range("B2").Value = "Ana Sophia Rocha"
$ws.Range("C2").Value = "Financeiro"
$ws.Range("D2").Value = "Consulta médica"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 45089
$ws.Range("G2").Value = 4631.22

# Row 3
$ws.Range("A3").Value = 93687
$ws.Range("B3").Value = "Erick da Rocha"
$ws.Range("C3").Value = "Recursos Humanos"
$ws.Range("D3").Value = "Outros"
$ws.Range("E3").Value = 4
$ws.Range("F3").Value = 45096
$ws.Range("G3").Value = 11506.59

# Row 4
$ws.Range("A4").Value = 18833
$ws.Range("B4").Value = "Nicole Teixeira"
$ws.Range("C4").Value = "Jurídico"
$ws.Range("D4").Value = "Consulta médica"
$ws.Range("E4").Value = 5
$ws.Range("F4").Value = 45103
$ws.Range("G4").Value = 12361.63

# Row 5
$ws.Range("A5").Value = 88740
$ws.Range("B5").Value = "Murilo Costela"
$ws.Range("C5").Value = "Financeiro"
$ws.Range("D5").Value = "Consulta médica"
$ws.Range("E5").Value = 5
$ws.Range("F5").Value = 45081
$ws.Range("G5").Value = 3670.02

# Row 6
$ws.Range("A6").Value = 72655
$ws.Range("B6").Value = "Sofia Nunes"
$ws.Range("C6").Value = "Operações"
$ws.Range("D6").Value = "Viagem de negócios"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 45081
$ws.Range("G6").Value = 5297.2

# Row 7
$ws.Range("A7").Value = 7289
$ws.Range("B7").Value = "Diego Moraes"
$ws.Range("C7").Value = "Recursos Humanos"
$ws.Range("D7").Value = "Consulta médica"
$ws.Range("E7").Value = 4
$ws.Range("F7").Value = 45094
$ws.Range("G7").Value = 9902.190000000001

# Row 8
$ws.Range("A8").Value = 36826
$ws.Range("B8").Value = "Arthur Correia"
$ws.Range("C8").Value = "Vendas"
$ws.Range("D8").Value = "Outros"
$ws.Range("E8").Value = 4
$ws.Range("F8").Value = 45087
$ws.Range("G8").Value = 6722.55

# Row 9
$ws.Range("A9").Value = 67508
$ws.Range("B9").Value = "Emanuella Viana"
$ws.Range("C9").Value = "Vendas"
$ws.Range("D9").Value = "Problemas pessoais"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 45103
$ws.Range("G9").Value = 5528.78

# Row 10
$ws.Range("A10").Value = 27792
$ws.Range("B10").Value = "Levi Gomes"
$ws.Range("C10").Value = "Recursos Humanos"
$ws.Range("D10").Value = "Consulta médica"
$ws.Range("E10").Value = 8
$ws.Range("F10").Value = 45097
$ws.Range("G10").Value = 5024.64

# Row 11
$ws.Range("A11").Value = 69106
$ws.Range("B11").Value = "Dra. Maria Clara da Rocha"
$ws.Range("C11").Value = "Financeiro"
$ws.Range("D11").Value = "Outros"
$ws.Range("E11").Value = 7
$ws.Range("F11").Value = 45095
$ws.Range("G11").Value = 10983.75
